# Applies the "Endless Symphony of Stars" -> "Marvelous World of Physics"
# rewrite described by the commit diff.
#
# Find.Execute signature used throughout:
#   Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#           MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#           Format, ReplaceWith, Replace)
# with Replace = 2 (wdReplaceAll).

$d = $word.ActiveDocument

# NOTE: Find.Execute's own ReplaceWith path silently runs the inserted text
# through "smart quotes" AutoFormat (straight ' / " become curly ' / ").
# Assigning the matched Range's .Text property does not, so we locate the
# text with Find (Replace = 0, i.e. wdReplaceNone) and then overwrite the
# Range's .Text directly to keep straight quotes/apostrophes intact.
function Replace-Text($old, $new, $wholeWord = $false) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $wholeWord, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
        return
    }
    $r.Text = $new
}

# --- Title / byline / contact block ---------------------------------------
Replace-Text "The Endless Symphony of Stars" "The Marvelous World of Physics: Unraveling the Mysteries of the Universe"
Replace-Text "Johnathan Orion" "Alex Rodriguez"
Replace-Text "johnathan" "alex"
Replace-Text "orion@starduststudies" "rodriguez@eduworld"
Replace-Text "com" "org" $true

# --- First body paragraph ---------------------------------------------------
# Sentence 1
Replace-Text "In the vast expanse of the universe, a majestic symphony of stars orchestrates a celestial performance of beauty and wonder" `
    "The realm of Physics presents us with an extraordinary odyssey into the very fabric of reality, where the laws governing the universe unfold, unveiling the intricacies of matter, energy, space, and time"

# Sentence 2
Replace-Text " Each star, a luminous celestial body, radiates with a brilliance unique to itself, forming intricate patterns and cosmic tapestries" `
    " Like a tapestry woven with enigmatic threads, Physics challenges our understanding, inviting us to unravel the mysteries of the cosmos"

# Sentence 3, plus two brand-new sentences appended right before the
# existing trailing "." run (so that old period now ends the new text).
Replace-Text " The night sky, with its myriad celestial dancers, has captured the imaginations of poets, astronomers, and stargazers throughout history" `
    " As we embark on this intellectual adventure, we will traverse the frontiers of knowledge, delving into the fundamental forces that shape our existence. From the symphony of particles to the enigmatic dance of quantum mechanics, Physics offers a mesmerizing journey into the unknown"

# Sentence 4 (first sentence after the double <w:br/>)
Replace-Text "These celestial orbs, dwelling in distant constellations, narrate captivating tales of life, evolution, and the formation of celestial bodies" `
    "As we journey through this enigmatic realm, we will probe the mysteries of the universe, delving into the depths of space and time"

# Sentence 5
Replace-Text " From the blazing fires of young, blue giants to the serene glow of red dwarfs, stars exemplify the diversity and grandeur of the cosmos" `
    " We will unravel the secrets of black holes, ponder the nature of dark matter, and explore the vastness of the cosmos"

# Sentence 6, plus two new sentences appended before the trailing "."
Replace-Text " Astronomers, with their telescopes pointed towards the celestial theater, analyze starlight, decipher its secrets, and unveil the mysteries of stellar physics" `
    " The journey of Physics will challenge our assumptions, expand our horizons, and inspire us with its profound beauty. Through this exploration, we will gain a deeper appreciation for the universe we inhabit and our place within it"

# Sentence 7 (first sentence after the second double <w:br/>)
Replace-Text "As conscientious observers, we are granted a privileged vantage point from which to appreciate the celestial symphony playing out before us" `
    "The quest for knowledge in Physics is an ever-evolving endeavor, constantly pushing the boundaries of human understanding"

# Sentence 8 (last sentence of the paragraph), plus two new sentences
# appended before the paragraph's final "."
Replace-Text " Whether it is gazing upon the Milky Way's ethereal luminescence or seeking out distant galaxies, observing stars enables us to probe the unfathomable depths of the cosmos, to traverse time and space through the medium of light" `
    " As we delve into the complexities of the physical world, we uncover new insights, unraveling the enigmas that have perplexed humanity for ages. With each discovery, we take another step towards comprehending the intricate workings of the universe, experiencing the exhilarating thrill of exploration and the profound satisfaction of unlocking nature's secrets"

# --- Summary paragraph -------------------------------------------------------
Replace-Text "The celestial symphony of stars captivates our imaginations with its beauty and grandeur, offering glimpses into the mysteries and wonders of the cosmos" `
    "Physics, with its intricate tapestry of theories and principles, unravels the enigmatic mysteries of the universe, inviting us to embark on a breathtaking odyssey of discovery"

Replace-Text " From fiery young giants to tranquil red dwarfs, stars radiate with diverse brilliance, forming cosmic tapestries that inspire poets and astronomers alike" `
    " This exploration of matter, energy, space, and time transcends the boundaries of human understanding, challenging our assumptions and inspiring awe with its profound beauty"

Replace-Text " Through the analysis of starlight, we unravel the intricate nature of stellar physics, bridging the gap between earthly and celestial realms" `
    " The journey through Physics pushes the limits of knowledge, revealing new insights into the cosmos, black holes, dark matter, and the boundless expanse of the universe"

Replace-Text " Stargazing becomes a profound experience, allowing us to connect with the vast universe and appreciate the symphony of light performed by these celestial wonders" `
    " Each step forward in Physics brings us closer to comprehending the intricate workings of nature, fulfilling humanity's insatiable quest for understanding the universe we inhabit"

# --- Trailing empty paragraph added at the end of the document -------------
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

Write-Output "done"
